$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF ("final delta S") column (F) for the rows that were
# recalculated after repulling the data / pushing all data / recomputing
# the mean.
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 5
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = -5
$ws.Range("F13").Value = -3
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = -3
